# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages rebuild at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of row -> [old value, new value] for column F, identical on both sheets
# except the last row index differs (11 on 展览, 15 on 全部类型).

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 8452
$ws1.Range("F5").Value = 6158
$ws1.Range("F6").Value = 532
$ws1.Range("F7").Value = 111
$ws1.Range("F10").Value = 321
$ws1.Range("F11").Value = 1123

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 8452
$ws4.Range("F5").Value = 6158
$ws4.Range("F6").Value = 532
$ws4.Range("F7").Value = 111
$ws4.Range("F10").Value = 321
$ws4.Range("F15").Value = 1123
